# Auto-generated script to append new survey response rows (452-463)
# matching the target diff for literacy_230925_tmp.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (451) down to the
# new rows (452-463) so the appended cells reuse the same cellXfs (s="2" for
# the timestamp column, s="3" for the rest) instead of creating new styles.
$ws.Range("A451:AD451").Copy()
$ws.Range("A452:AD463").PasteSpecial(-4122)
$excel.CutCopyMode = 0


# Row 452
$ws.Cells.Item(452, 1).Value = 45200.9405944676
$ws.Cells.Item(452, 2).Value = "scw0922@naver.com"
$ws.Cells.Item(452, 3).Value = "간호학과"
$ws.Cells.Item(452, 4).Value = 20236256
$ws.Cells.Item(452, 5).Value = "신채원"
$ws.Cells.Item(452, 6).Value = 3
$ws.Cells.Item(452, 7).Value = "2. 시세 정보는 약 1개월 간격으로 갱신된다."
$ws.Cells.Item(452, 8).Value = "2. 인공적인 향기가 인체에 해롭지 않을까요?"
$ws.Cells.Item(452, 9).Value = "4. 6630번 버스를 타면 한마음예식장에 갈 수 있다."
$ws.Cells.Item(452, 10).Value = "3. 거실 바닥을 자주 물걸레로 닦는다"
$ws.Cells.Item(452, 11).Value = "3. 음주 운전이 의심될 경우 경찰관은 바로 운전자에게 혈액 채취를 명할 수 있군."
$ws.Cells.Item(452, 12).Value = "2. 친구를 만나서 가까운 산에 오른다."
$ws.Cells.Item(452, 13).Value = "3. 허위로 신고하면 10만 원의 과태료를 물게 된다."
$ws.Cells.Item(452, 14).Value = "4. 8,000 원"
$ws.Cells.Item(452, 15).Value = "3. 기침 감기에 종합 감기약을 먹으면 기침약을 먹은 것과 효과가 같다."
$ws.Cells.Item(452, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(452, 17).Value = "3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요."
$ws.Cells.Item(452, 18).Value = "3. 홍길동 씨가 보내려는 돈은 30,500 원이다."
$ws.Cells.Item(452, 19).Value = "2. 보증 기간 동안에는 건전지를 무상으로 제공한다."
$ws.Cells.Item(452, 20).Value = "2. 동남쪽에서부터 꽃이 피기 시작한다."
$ws.Cells.Item(452, 21).Value = "2. 벽지를 구입한 고객에게는 대걸레를 준다."
$ws.Cells.Item(452, 22).Value = "2. 오전 6시"
$ws.Cells.Item(452, 23).Value = "1. 내 전공이 화학이니 지원 가능하겠군."
$ws.Cells.Item(452, 24).Value = "3. 상점 주소"
$ws.Cells.Item(452, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(452, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(452, 27).Value = "2. 교육은 특정 요일에 실시된다."
$ws.Cells.Item(452, 28).Value = "3. 공공장소에서는 전화 예절을 지켜야 한다."
$ws.Cells.Item(452, 29).Value = "4. 16권의 만화 '토지'가 한꺼번에 출간되었다."
$ws.Cells.Item(452, 30).Value = "4. 야구 중계는 오후 2시 25분에 시작한다."

# Row 453
$ws.Cells.Item(453, 1).Value = 45200.94394950231
$ws.Cells.Item(453, 2).Value = "harin3040@naver.com"
$ws.Cells.Item(453, 3).Value = "심리학과"
$ws.Cells.Item(453, 4).Value = 20232113
$ws.Cells.Item(453, 5).Value = "김현진"
$ws.Cells.Item(453, 6).Value = 3
$ws.Cells.Item(453, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(453, 8).Value = "4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?"
$ws.Cells.Item(453, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(453, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(453, 11).Value = "2. 운전자의 음주 운전 여부에 대한 최종 판단은 혈액 채취 결과만 인정이 되는군."
$ws.Cells.Item(453, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(453, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(453, 14).Value = "3. 7,000 원"
$ws.Cells.Item(453, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(453, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(453, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(453, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(453, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(453, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(453, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(453, 22).Value = "4. 오후3시"
$ws.Cells.Item(453, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(453, 24).Value = "1. 상품 가격"
$ws.Cells.Item(453, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(453, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(453, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(453, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(453, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(453, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 454
$ws.Cells.Item(454, 1).Value = 45200.94654295139
$ws.Cells.Item(454, 2).Value = "shanesun0923@gmail.com"
$ws.Cells.Item(454, 3).Value = "간호학과"
$ws.Cells.Item(454, 4).Value = 20236253
$ws.Cells.Item(454, 5).Value = "선세인"
$ws.Cells.Item(454, 6).Value = 3
$ws.Cells.Item(454, 7).Value = "3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다."
$ws.Cells.Item(454, 8).Value = "4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?"
$ws.Cells.Item(454, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(454, 10).Value = "3. 거실 바닥을 자주 물걸레로 닦는다"
$ws.Cells.Item(454, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(454, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(454, 13).Value = "1. 이 서식은 국내 전입신고 시에만 사용할 수 있다."
$ws.Cells.Item(454, 14).Value = "3. 7,000 원"
$ws.Cells.Item(454, 15).Value = "4. 남은 약은 반드시 냉장고에 보관해야 한다."
$ws.Cells.Item(454, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(454, 17).Value = "4. 확인증을 잃어버렸는데, 다시 발급해 주겠지."
$ws.Cells.Item(454, 18).Value = "3. 홍길동 씨가 보내려는 돈은 30,500 원이다."
$ws.Cells.Item(454, 19).Value = "2. 보증 기간 동안에는 건전지를 무상으로 제공한다."
$ws.Cells.Item(454, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(454, 21).Value = "4. 조기 품절 시에는 사은품이 다른 물품으로 대체될 수 있다."
$ws.Cells.Item(454, 22).Value = "3. 오후 6시"
$ws.Cells.Item(454, 23).Value = "4. 일주일에 이틀은 쉴 수 있겠군."
$ws.Cells.Item(454, 24).Value = "3. 상점 주소"
$ws.Cells.Item(454, 25).Value = "1. 뜻풀이 '1'"
$ws.Cells.Item(454, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(454, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(454, 28).Value = "4. 공공장소에서는 떠들지 말아야 한다."
$ws.Cells.Item(454, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(454, 30).Value = "4. 야구 중계는 오후 2시 25분에 시작한다."

# Row 455
$ws.Cells.Item(455, 1).Value = 45200.95309476852
$ws.Cells.Item(455, 2).Value = "sung93716@gmail.com"
$ws.Cells.Item(455, 3).Value = "데이터사이언스학부"
$ws.Cells.Item(455, 4).Value = 20233261
$ws.Cells.Item(455, 5).Value = "한예림"
$ws.Cells.Item(455, 6).Value = 3
$ws.Cells.Item(455, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(455, 8).Value = "1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?"
$ws.Cells.Item(455, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(455, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(455, 11).Value = "2. 운전자의 음주 운전 여부에 대한 최종 판단은 혈액 채취 결과만 인정이 되는군."
$ws.Cells.Item(455, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(455, 13).Value = "3. 허위로 신고하면 10만 원의 과태료를 물게 된다."
$ws.Cells.Item(455, 14).Value = "2. 6,000 원"
$ws.Cells.Item(455, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(455, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(455, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(455, 18).Value = "4. 홍길동 씨는 세계은행에서 송금을 하고 있다."
$ws.Cells.Item(455, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(455, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(455, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(455, 22).Value = "4. 오후3시"
$ws.Cells.Item(455, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(455, 24).Value = "1. 상품 가격"
$ws.Cells.Item(455, 25).Value = "2. 뜻풀이 '2'"
$ws.Cells.Item(455, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(455, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(455, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(455, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(455, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 456
$ws.Cells.Item(456, 1).Value = 45200.968275844905
$ws.Cells.Item(456, 2).Value = "dksdksqh1018@naver.com"
$ws.Cells.Item(456, 3).Value = "미디어스쿨"
$ws.Cells.Item(456, 4).Value = 20232549
$ws.Cells.Item(456, 5).Value = "안보민"
$ws.Cells.Item(456, 6).Value = 3
$ws.Cells.Item(456, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(456, 8).Value = "4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?"
$ws.Cells.Item(456, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(456, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(456, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(456, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(456, 13).Value = "3. 허위로 신고하면 10만 원의 과태료를 물게 된다."
$ws.Cells.Item(456, 14).Value = "2. 6,000 원"
$ws.Cells.Item(456, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(456, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(456, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(456, 18).Value = "4. 홍길동 씨는 세계은행에서 송금을 하고 있다."
$ws.Cells.Item(456, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(456, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(456, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(456, 22).Value = "4. 오후3시"
$ws.Cells.Item(456, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(456, 24).Value = "1. 상품 가격"
$ws.Cells.Item(456, 25).Value = "1. 뜻풀이 '1'"
$ws.Cells.Item(456, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(456, 27).Value = "2. 교육은 특정 요일에 실시된다."
$ws.Cells.Item(456, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(456, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(456, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 457
$ws.Cells.Item(457, 1).Value = 45200.971810937495
$ws.Cells.Item(457, 2).Value = "gaejisub@gmail.com"
$ws.Cells.Item(457, 3).Value = "콘텐츠it"
$ws.Cells.Item(457, 4).Value = 20225169
$ws.Cells.Item(457, 5).Value = "배승유"
$ws.Cells.Item(457, 6).Value = 3
$ws.Cells.Item(457, 7).Value = "3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다."
$ws.Cells.Item(457, 8).Value = "1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?"
$ws.Cells.Item(457, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(457, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(457, 11).Value = "3. 음주 운전이 의심될 경우 경찰관은 바로 운전자에게 혈액 채취를 명할 수 있군."
$ws.Cells.Item(457, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(457, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(457, 14).Value = "2. 6,000 원"
$ws.Cells.Item(457, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(457, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(457, 17).Value = "3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요."
$ws.Cells.Item(457, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(457, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(457, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(457, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(457, 22).Value = "4. 오후3시"
$ws.Cells.Item(457, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(457, 24).Value = "1. 상품 가격"
$ws.Cells.Item(457, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(457, 26).Value = "3. 여우비, 소나무향기"
$ws.Cells.Item(457, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(457, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(457, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(457, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 458
$ws.Cells.Item(458, 1).Value = 45200.981516238426
$ws.Cells.Item(458, 2).Value = "gustj1654@naver.com"
$ws.Cells.Item(458, 3).Value = "심리학과"
$ws.Cells.Item(458, 4).Value = 20232137
$ws.Cells.Item(458, 5).Value = "조현서"
$ws.Cells.Item(458, 6).Value = 3
$ws.Cells.Item(458, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(458, 8).Value = "1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?"
$ws.Cells.Item(458, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(458, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(458, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(458, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(458, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(458, 14).Value = "2. 6,000 원"
$ws.Cells.Item(458, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(458, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(458, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(458, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(458, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(458, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(458, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(458, 22).Value = "4. 오후3시"
$ws.Cells.Item(458, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(458, 24).Value = "1. 상품 가격"
$ws.Cells.Item(458, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(458, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(458, 27).Value = "1. 이 프로그램은 노인들만을 위한 것이다."
$ws.Cells.Item(458, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(458, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(458, 30).Value = "3. 'TV는 사랑을 싣고'는 다시 보기를 제공하지 않는다."

# Row 459
$ws.Cells.Item(459, 1).Value = 45200.986926377314
$ws.Cells.Item(459, 2).Value = "yejin4259@naver.com"
$ws.Cells.Item(459, 3).Value = "언어청각학부"
$ws.Cells.Item(459, 4).Value = 20233951
$ws.Cells.Item(459, 5).Value = "이예진"
$ws.Cells.Item(459, 6).Value = 2
$ws.Cells.Item(459, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(459, 8).Value = "4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?"
$ws.Cells.Item(459, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(459, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(459, 11).Value = "2. 운전자의 음주 운전 여부에 대한 최종 판단은 혈액 채취 결과만 인정이 되는군."
$ws.Cells.Item(459, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(459, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(459, 14).Value = "1. 5,000 원"
$ws.Cells.Item(459, 15).Value = "4. 남은 약은 반드시 냉장고에 보관해야 한다."
$ws.Cells.Item(459, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(459, 17).Value = "3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요."
$ws.Cells.Item(459, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(459, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(459, 20).Value = "2. 동남쪽에서부터 꽃이 피기 시작한다."
$ws.Cells.Item(459, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(459, 22).Value = "1. 오전 3시"
$ws.Cells.Item(459, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(459, 24).Value = "1. 상품 가격"
$ws.Cells.Item(459, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(459, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(459, 27).Value = "2. 교육은 특정 요일에 실시된다."
$ws.Cells.Item(459, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(459, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(459, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 460
$ws.Cells.Item(460, 1).Value = 45200.9901277662
$ws.Cells.Item(460, 2).Value = "ghys1837@naver.com"
$ws.Cells.Item(460, 3).Value = "언어청각학부 청각학 전공"
$ws.Cells.Item(460, 4).Value = 20213939
$ws.Cells.Item(460, 5).Value = "안영서"
$ws.Cells.Item(460, 6).Value = 3
$ws.Cells.Item(460, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(460, 8).Value = "2. 인공적인 향기가 인체에 해롭지 않을까요?"
$ws.Cells.Item(460, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(460, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(460, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(460, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(460, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(460, 14).Value = "2. 6,000 원"
$ws.Cells.Item(460, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(460, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(460, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(460, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(460, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(460, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(460, 21).Value = "2. 벽지를 구입한 고객에게는 대걸레를 준다."
$ws.Cells.Item(460, 22).Value = "1. 오전 3시"
$ws.Cells.Item(460, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(460, 24).Value = "1. 상품 가격"
$ws.Cells.Item(460, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(460, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(460, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(460, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(460, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(460, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 461
$ws.Cells.Item(461, 1).Value = 45200.99886603009
$ws.Cells.Item(461, 2).Value = "tjdus3641@gmail.com"
$ws.Cells.Item(461, 3).Value = "간호학과"
$ws.Cells.Item(461, 4).Value = 20226283
$ws.Cells.Item(461, 5).Value = "장서연"
$ws.Cells.Item(461, 6).Value = 3
$ws.Cells.Item(461, 7).Value = "3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다."
$ws.Cells.Item(461, 8).Value = "4. 향기를 전달할 수 있는 휴대전화의 가격은 얼마입니까?"
$ws.Cells.Item(461, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(461, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(461, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(461, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(461, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(461, 14).Value = "2. 6,000 원"
$ws.Cells.Item(461, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(461, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(461, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(461, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(461, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(461, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(461, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(461, 22).Value = "4. 오후3시"
$ws.Cells.Item(461, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(461, 24).Value = "1. 상품 가격"
$ws.Cells.Item(461, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(461, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(461, 27).Value = "2. 교육은 특정 요일에 실시된다."
$ws.Cells.Item(461, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(461, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(461, 30).Value = "4. 야구 중계는 오후 2시 25분에 시작한다."

# Row 462
$ws.Cells.Item(462, 1).Value = 45201.002475104164
$ws.Cells.Item(462, 2).Value = "rkqls3333@gmail.com"
$ws.Cells.Item(462, 3).Value = "간호학과"
$ws.Cells.Item(462, 4).Value = 20236205
$ws.Cells.Item(462, 5).Value = "권가빈"
$ws.Cells.Item(462, 6).Value = 3
$ws.Cells.Item(462, 7).Value = "3. 3월에서 6월까지 매매 가격이나 전세 가격 모두 변화가 없다."
$ws.Cells.Item(462, 8).Value = "1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?"
$ws.Cells.Item(462, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(462, 10).Value = "3. 거실 바닥을 자주 물걸레로 닦는다"
$ws.Cells.Item(462, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(462, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(462, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(462, 14).Value = "2. 6,000 원"
$ws.Cells.Item(462, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(462, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(462, 17).Value = "2. 오늘은 4월 30일이니까, 국립중앙박물관에서 사용할 수 있겠죠?"
$ws.Cells.Item(462, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(462, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(462, 20).Value = "3. 벚꽃이 개나리보다 일주일 정도 일찍 핀다."
$ws.Cells.Item(462, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(462, 22).Value = "1. 오전 3시"
$ws.Cells.Item(462, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(462, 24).Value = "1. 상품 가격"
$ws.Cells.Item(462, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(462, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(462, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(462, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(462, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(462, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

# Row 463
$ws.Cells.Item(463, 1).Value = 45201.00958046297
$ws.Cells.Item(463, 2).Value = "nanadiana222@naver.com"
$ws.Cells.Item(463, 3).Value = "빅데이터과"
$ws.Cells.Item(463, 4).Value = 202252161
$ws.Cells.Item(463, 5).Value = "조희진"
$ws.Cells.Item(463, 6).Value = 3
$ws.Cells.Item(463, 7).Value = "4. 2월과 3월 사이에 매매 가격이 1억 원 정도 상승했다."
$ws.Cells.Item(463, 8).Value = "1. 전달 가능한 향기가 겨우 16가지밖에 안 됩니까?"
$ws.Cells.Item(463, 9).Value = "3. 한마음예식장은 주차 공간이 충분하므로 승용차를 이용하는 것이 좋다."
$ws.Cells.Item(463, 10).Value = "2. 채소 등 농산물은 익혀 먹는다"
$ws.Cells.Item(463, 11).Value = "4. 혈중 알코올 농도가 0.05퍼센트 미만이면 운전을 못할 정도로 취한 상태는 아니라고 보는군."
$ws.Cells.Item(463, 12).Value = "1. 학교 가는 아이에게 우산을 챙겨 준다."
$ws.Cells.Item(463, 13).Value = "4. 정당한 이유가 없는 한 14일 이내에 신고해야 한다."
$ws.Cells.Item(463, 14).Value = "2. 6,000 원"
$ws.Cells.Item(463, 15).Value = "2. 약은 우유나 주스와 함께 먹지 않는 것이 좋다."
$ws.Cells.Item(463, 16).Value = "1. 이 공지 사항은 '봉투나라'의 관리자가 작성한 것이다."
$ws.Cells.Item(463, 17).Value = "3. 이 주차장은 할인이 안 될지 모르니 미리 물어봐야 해요."
$ws.Cells.Item(463, 18).Value = "2. 홍길동 씨는 김명자 씨에게 송금을 하고 있다."
$ws.Cells.Item(463, 19).Value = "4. 수리할 수 없는 제품은 새것으로 교환해 준다."
$ws.Cells.Item(463, 20).Value = "2. 동남쪽에서부터 꽃이 피기 시작한다."
$ws.Cells.Item(463, 21).Value = "1. 무료 배송을 받으려면 5만 원어치 이상을 사야 된다."
$ws.Cells.Item(463, 22).Value = "4. 오후3시"
$ws.Cells.Item(463, 23).Value = "2. 남녀를 차별 대우하다니 한심하네."
$ws.Cells.Item(463, 24).Value = "1. 상품 가격"
$ws.Cells.Item(463, 25).Value = "3. 뜻풀이 '3'"
$ws.Cells.Item(463, 26).Value = "2. 여우비, 장난감박물관"
$ws.Cells.Item(463, 27).Value = "4. 5만원을 내면 모든 강의를 들을 수 있다."
$ws.Cells.Item(463, 28).Value = "1. 출산율을 높여야 한다."
$ws.Cells.Item(463, 29).Value = "3. '토지'에는 700여 명의 인물이 등장한다."
$ws.Cells.Item(463, 30).Value = "2. 낮 12시 뉴스는 자막 방송을 하지 않는다."

Write-Host "Appended rows 452-463"